# Updated code with StudyId
# Target sheet: "TS"
# 1) Remove TS-parameter rows that have no mapped value (blank TSVAL in column G),
#    compacting the remaining mapped rows upward.
# 2) Stamp the STUDYID (column A) as "CB0321" on every remaining data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TS")

# Rows (1-based, as in the original sheet) whose TSVAL (column G) is blank and
# therefore should be dropped. Deleted from the bottom up so earlier row
# numbers in this list stay valid as the sheet shrinks.
$rowsToDelete = @(60,58,56,55,52,51,50,48,47,46,45,44,43,41,40,39,36,35,34,29,28,27,25,24,20,19,18,17,12,11,9,8,7,6,5)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).EntireRow.Delete()
}

# After compaction, the kept TS-parameter rows occupy rows 2 through 25.
$lastRow = 25
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "CB0321"
}
